$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new header row above the existing data, shifting rows 1-8 down to 2-9 ---
$ws.Rows.Item(1).Insert()

# Copy the look of an existing data-style cell (wrap text, Georgia 10pt) onto the header row
$ws.Range("B2").Copy()
$hdr = $ws.Range("A1:J1")
$hdr.PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# New header captions
$headers = @(
    "Name of organisation",
    "Description of organisation",
    "Mission/ Objectives/ Purpose",
    "Programmes/ projects",
    "Funding sources",
    "Collaboration with government / businesses",
    "Choice of Climate action",
    "No. of employees",
    "Geographical focus",
    "Nationality"
)
for ($c = 1; $c -le 10; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# Shade the header row
$hdr.Interior.Color = 13421772
$ws.Rows.Item(1).RowHeight = 57

# Snapshot the organisation-name column's current look (Georgia 10pt underlined blue,
# style already carried down by the row insert) so it can be restored after adding hyperlinks.
$ws.Range("Z1").Clear()
$ws.Range("A2").Copy($ws.Range("Z1"))
$excel.CutCopyMode = $false

# --- Re-create the hyperlinks on the shifted organisation-name column (now A2:A9) ---
$urls = @(
    "http://mlup-baitong.org/",
    "https://www.savetheearthcambodia.org/",
    "https://songsaafoundation.org/",
    "https://www.ccc-cambodia.org/en/ngodb/ngo-information/4441",
    "http://cepa-cambodia.org/",
    "https://centerforsustainablewater.org/",
    "https://cambodia.wcs.org/",
    "https://savetheearthinternational.org/"
)
for ($i = 0; $i -lt 8; $i++) {
    $ws.Hyperlinks.Add($ws.Range("A$($i + 2)"), $urls[$i]) | Out-Null
}

# Restore the original (non built-in-hyperlink) look that Hyperlinks.Add overwrote
$ws.Range("Z1").Copy()
$ws.Range("A2:A9").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("Z1").Clear()

# Selection now rests on the header row only
$ws.Range("A1:J1").Select()
